$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 812.4286
$ws.Range("I19").Value = 659
$ws.Range("J19").Value = 927.5
$ws.Range("K19").Value = 659
$ws.Range("L19").Value = 927.5
$ws.Range("M19").Value = -484
$ws.Range("N19").Value = -1277.5
$ws.Range("H21").Value = 26652.334
$ws.Range("J21").Value = 26652.334
$ws.Range("L21").Value = 26652.334
$ws.Range("N21").Value = -27588.334
$ws.Range("H23").Value = 26652.334
$ws.Range("J23").Value = 26652.334
$ws.Range("L23").Value = 26652.334
$ws.Range("N23").Value = -27120.334
$ws.Range("H33").Value = 380.4
$ws.Range("J33").Value = 522.6
$ws.Range("L33").Value = 522.6
$ws.Range("N33").Value = -980.6
$ws.Range("H46").Value = 250
$ws.Range("I46").Value = 250
$ws.Range("K46").Value = 750
$ws.Range("M46").Value = -631
$ws.Range("H60").Value = 250
$ws.Range("I60").Value = 250
$ws.Range("K60").Value = 750
$ws.Range("M60").Value = -266
$ws.Range("H98").Value = 1819.4918
$ws.Range("I98").Value = 1929.1455
$ws.Range("K98").Value = 1929.1455
$ws.Range("M98").Value = -431.1455000000001
$ws.Range("H107").Value = 4362
$ws.Range("I107").Value = 6966.6665
$ws.Range("J107").Value = 2799.2
$ws.Range("K107").Value = 6966.6665
$ws.Range("L107").Value = 2799.2
$ws.Range("M107").Value = -5046.6665
$ws.Range("N107").Value = -6639.2
$ws.Range("H113").Value = 3902.5
$ws.Range("I113").Value = 4005
$ws.Range("K113").Value = 4005
$ws.Range("M113").Value = -751
$ws.Range("H116").Value = 2839.4
$ws.Range("I116").Value = 2065
$ws.Range("J116").Value = 4001
$ws.Range("K116").Value = 2065
$ws.Range("L116").Value = 4001
$ws.Range("M116").Value = 1377
$ws.Range("N116").Value = -10885
$ws.Range("H122").Value = 1819.4918
$ws.Range("I122").Value = 1929.1455
$ws.Range("K122").Value = 5787.4365
$ws.Range("M122").Value = -3337.4365
$ws.Range("H132").Value = 11119100
$ws.Range("I132").Value = 12351788
$ws.Range("J132").Value = 24901.666
$ws.Range("K132").Value = 37055364
$ws.Range("L132").Value = 74704.99800000001
$ws.Range("M132").Value = -37052834
$ws.Range("N132").Value = -79764.99800000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11780.182
$ws.Range("I32").Value = 10100.704
$ws.Range("K32").Value = 10100.704
$ws.Range("M32").Value = -9813.704
$ws.Range("H61").Value = 100001970
$ws.Range("I61").Value = 166667950
$ws.Range("K61").Value = 166667950
$ws.Range("M61").Value = -166667738
$ws.Range("H74").Value = 1739.3334
$ws.Range("I74").Value = 1141.3334
$ws.Range("K74").Value = 1141.3334
$ws.Range("M74").Value = -267.3334
$ws.Range("H77").Value = 1739.3334
$ws.Range("I77").Value = 1141.3334
$ws.Range("K77").Value = 5706.666999999999
$ws.Range("M77").Value = -1338.666999999999
$ws.Range("H110").Value = 2226.7
$ws.Range("I110").Value = 1373.5
$ws.Range("J110").Value = 3506.5
$ws.Range("K110").Value = 1373.5
$ws.Range("L110").Value = 3506.5
$ws.Range("M110").Value = 671.5
$ws.Range("N110").Value = -7596.5
$ws.Range("H136").Value = 100001970
$ws.Range("I136").Value = 166667950
$ws.Range("K136").Value = 500003850
$ws.Range("M136").Value = -500001300

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 731.4167
$ws.Range("I80").Value = 130.83333
$ws.Range("J80").Value = 1332
$ws.Range("K80").Value = 130.83333
$ws.Range("L80").Value = 1332
$ws.Range("M80").Value = 867.1666700000001
$ws.Range("N80").Value = -3328
$ws.Range("H83").Value = 731.4167
$ws.Range("I83").Value = 130.83333
$ws.Range("J83").Value = 1332
$ws.Range("K83").Value = 654.1666499999999
$ws.Range("L83").Value = 6660
$ws.Range("M83").Value = 4337.83335
$ws.Range("N83").Value = -16644
$ws.Range("H94").Value = 20834002
$ws.Range("I94").Value = 25000302
$ws.Range("J94").Value = 2505
$ws.Range("K94").Value = 25000302
$ws.Range("L94").Value = 2505
$ws.Range("M94").Value = -24999851
$ws.Range("N94").Value = -3407
$ws.Range("H105").Value = 90910400
$ws.Range("I105").Value = 100001140
$ws.Range("K105").Value = 100001140
$ws.Range("M105").Value = -99999393
$ws.Range("H107").Value = 3108
$ws.Range("I107").Value = 2011
$ws.Range("J107").Value = 3656.5
$ws.Range("K107").Value = 2011
$ws.Range("L107").Value = 3656.5
$ws.Range("M107").Value = -91
$ws.Range("N107").Value = -7496.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1240.5
$ws.Range("I31").Value = 1022.375
$ws.Range("J31").Value = 3421.75
$ws.Range("K31").Value = 1022.375
$ws.Range("L31").Value = 3421.75
$ws.Range("M31").Value = -727.375
$ws.Range("N31").Value = -4011.75
$ws.Range("H34").Value = 1240.5
$ws.Range("I34").Value = 1022.375
$ws.Range("J34").Value = 3421.75
$ws.Range("K34").Value = 1022.375
$ws.Range("L34").Value = 3421.75
$ws.Range("M34").Value = -820.375
$ws.Range("N34").Value = -3825.75
$ws.Range("H141").Value = 619765.3
$ws.Range("J141").Value = 619765.3
$ws.Range("L141").Value = 619765.3
$ws.Range("N141").Value = -630125.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1628
$ws.Range("I34").Value = 599
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1797
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -1713
$ws.Range("N34").Value = -9168
$ws.Range("H92").Value = 266.46155
$ws.Range("I92").Value = 246.94118
$ws.Range("J92").Value = 303.33334
$ws.Range("K92").Value = 740.82354
$ws.Range("L92").Value = 910.0000200000001
$ws.Range("M92").Value = 507.17646
$ws.Range("N92").Value = -3406.00002
$ws.Range("H107").Value = 6568.3125
$ws.Range("J107").Value = 8589.5
$ws.Range("L107").Value = 25768.5
$ws.Range("N107").Value = -29608.5
$ws.Range("H122").Value = 1067.3636
$ws.Range("I122").Value = 828
$ws.Range("K122").Value = 7452
$ws.Range("M122").Value = -5002
$ws.Range("H131").Value = 38467200
$ws.Range("J131").Value = 6636.636
$ws.Range("L131").Value = 19909.908
$ws.Range("N131").Value = -29989.908

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2509.2778
$ws.Range("I80").Value = 1544.375
$ws.Range("J80").Value = 3281.2
$ws.Range("K80").Value = 1544.375
$ws.Range("L80").Value = 3281.2
$ws.Range("M80").Value = -546.375
$ws.Range("N80").Value = -5277.2
$ws.Range("H83").Value = 2509.2778
$ws.Range("I83").Value = 1544.375
$ws.Range("J83").Value = 3281.2
$ws.Range("K83").Value = 7721.875
$ws.Range("L83").Value = 16406
$ws.Range("M83").Value = -2729.875
$ws.Range("N83").Value = -26390
$ws.Range("H126").Value = 2398.75
$ws.Range("I126").Value = 1850
$ws.Range("K126").Value = 5550
$ws.Range("M126").Value = -3080

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1085.1
$ws.Range("I55").Value = 1035.8572
$ws.Range("J55").Value = 1200
$ws.Range("K55").Value = 1035.8572
$ws.Range("L55").Value = 1200
$ws.Range("M55").Value = -862.8571999999999
$ws.Range("N55").Value = -1546

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2799.9285
$ws.Range("I96").Value = 1809.5
$ws.Range("K96").Value = 1809.5
$ws.Range("M96").Value = -436.5
$ws.Range("H122").Value = 12504618
$ws.Range("I122").Value = 16671992
$ws.Range("J122").Value = 2497
$ws.Range("K122").Value = 50015976
$ws.Range("L122").Value = 7491
$ws.Range("M122").Value = -50013526
$ws.Range("N122").Value = -12391
